$wb = $excel.ActiveWorkbook

# ALC row 4 (Leve Item ID 5470)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 13).ClearContents()

# ALC row 69 (Leve Item ID 12616)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 3730
$ws.Cells.Item(69, 9).Value = 3000
$ws.Cells.Item(69, 10).Value = 3912.5
$ws.Cells.Item(69, 11).Value = 9000
$ws.Cells.Item(69, 12).Value = 11737.5
$ws.Cells.Item(69, 13).Value = -8126
$ws.Cells.Item(69, 14).Value = -13485.5

# ALC row 72 (Leve Item ID 12616)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(72, 8).Value = 3730
$ws.Cells.Item(72, 9).Value = 3000
$ws.Cells.Item(72, 10).Value = 3912.5
$ws.Cells.Item(72, 11).Value = 27000
$ws.Cells.Item(72, 12).Value = 35212.5
$ws.Cells.Item(72, 13).Value = -22632
$ws.Cells.Item(72, 14).Value = -43948.5

# ALC row 112 (Leve Item ID 27960)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 3694.0625
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 10).Value = 3694.0625
$ws.Cells.Item(112, 11).Value = 0
$ws.Cells.Item(112, 12).Value = 11082.1875
$ws.Cells.Item(112, 13).ClearContents()
$ws.Cells.Item(112, 14).Value = -13298.1875

# ARM row 110 (Leve Item ID 27708)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 4558
$ws.Cells.Item(110, 9).Value = 4558
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 11).Value = 4558
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 13).Value = -2513
$ws.Cells.Item(110, 14).ClearContents()

# ARM row 132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1920.0526
$ws.Cells.Item(132, 9).Value = 1738.0385
$ws.Cells.Item(132, 10).Value = 2314.4167
$ws.Cells.Item(132, 11).Value = 5214.1155
$ws.Cells.Item(132, 12).Value = 6943.250100000001
$ws.Cells.Item(132, 13).Value = -2684.1155
$ws.Cells.Item(132, 14).Value = -12003.2501

# BSM row 53 (Leve Item ID 27158)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(53, 8).Value = 45020
$ws.Cells.Item(53, 10).Value = 45020
$ws.Cells.Item(53, 12).Value = 45020
$ws.Cells.Item(53, 14).Value = -46168

# BSM row 107 (Leve Item ID 27706)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 2353
$ws.Cells.Item(107, 9).Value = 1745.6666
$ws.Cells.Item(107, 11).Value = 1745.6666
$ws.Cells.Item(107, 13).Value = 174.3334

# BSM row 134 (Leve Item ID 43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1130.9565
$ws.Cells.Item(134, 9).Value = 1060.6
$ws.Cells.Item(134, 11).Value = 3181.8
$ws.Cells.Item(134, 13).Value = -646.7999999999997

# CRP row 16 (Leve Item ID 27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2494.182
$ws.Cells.Item(16, 9).Value = 2849
$ws.Cells.Item(16, 10).Value = 897.5
$ws.Cells.Item(16, 11).Value = 2849
$ws.Cells.Item(16, 12).Value = 897.5
$ws.Cells.Item(16, 13).Value = -2562
$ws.Cells.Item(16, 14).Value = -1471.5

# CRP row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1976.4637
$ws.Cells.Item(31, 9).Value = 1512.561
$ws.Cells.Item(31, 10).Value = 2655.75
$ws.Cells.Item(31, 11).Value = 1512.561
$ws.Cells.Item(31, 12).Value = 2655.75
$ws.Cells.Item(31, 13).Value = -1217.561
$ws.Cells.Item(31, 14).Value = -3245.75

# CRP row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 1976.4637
$ws.Cells.Item(34, 9).Value = 1512.561
$ws.Cells.Item(34, 10).Value = 2655.75
$ws.Cells.Item(34, 11).Value = 1512.561
$ws.Cells.Item(34, 12).Value = 2655.75
$ws.Cells.Item(34, 13).Value = -1310.561
$ws.Cells.Item(34, 14).Value = -3059.75

# CRP row 62 (Leve Item ID 12580)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 3895.5881
$ws.Cells.Item(62, 9).Value = 3808.3333
$ws.Cells.Item(62, 10).Value = 3993.75
$ws.Cells.Item(62, 11).Value = 3808.3333
$ws.Cells.Item(62, 12).Value = 3993.75
$ws.Cells.Item(62, 13).Value = -3184.3333
$ws.Cells.Item(62, 14).Value = -5241.75

# CRP row 65 (Leve Item ID 12580)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 3895.5881
$ws.Cells.Item(65, 9).Value = 3808.3333
$ws.Cells.Item(65, 10).Value = 3993.75
$ws.Cells.Item(65, 11).Value = 19041.6665
$ws.Cells.Item(65, 12).Value = 19968.75
$ws.Cells.Item(65, 13).Value = -15921.6665
$ws.Cells.Item(65, 14).Value = -26208.75

# CRP row 105 (Leve Item ID 19928)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 2035
$ws.Cells.Item(105, 9).Value = 1060
$ws.Cells.Item(105, 10).Value = 4960
$ws.Cells.Item(105, 11).Value = 1060
$ws.Cells.Item(105, 12).Value = 4960
$ws.Cells.Item(105, 13).Value = 687
$ws.Cells.Item(105, 14).Value = -8454

# CRP row 113 (Leve Item ID 27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 2494.182
$ws.Cells.Item(113, 9).Value = 2849
$ws.Cells.Item(113, 10).Value = 897.5
$ws.Cells.Item(113, 11).Value = 2849
$ws.Cells.Item(113, 12).Value = 897.5
$ws.Cells.Item(113, 13).Value = -679
$ws.Cells.Item(113, 14).Value = -5237.5

# CUL row 23 (Leve Item ID 4858)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 29.166666
$ws.Cells.Item(23, 10).Value = 43.333332
$ws.Cells.Item(23, 12).Value = 129.999996
$ws.Cells.Item(23, 14).Value = -599.999996

# CUL row 69 (Leve Item ID 12850)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 3560
$ws.Cells.Item(69, 9).Value = 1250
$ws.Cells.Item(69, 10).Value = 5100
$ws.Cells.Item(69, 11).Value = 3750
$ws.Cells.Item(69, 12).Value = 15300
$ws.Cells.Item(69, 13).Value = -2939
$ws.Cells.Item(69, 14).Value = -16922

# CUL row 72 (Leve Item ID 12850)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(72, 8).Value = 3560
$ws.Cells.Item(72, 9).Value = 1250
$ws.Cells.Item(72, 10).Value = 5100
$ws.Cells.Item(72, 11).Value = 11250
$ws.Cells.Item(72, 12).Value = 45900
$ws.Cells.Item(72, 13).Value = -7194
$ws.Cells.Item(72, 14).Value = -54012

# CUL row 75 (Leve Item ID 12863)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(75, 8).Value = 12000
$ws.Cells.Item(75, 10).Value = 12000
$ws.Cells.Item(75, 12).Value = 36000
$ws.Cells.Item(75, 14).Value = -37996

# CUL row 78 (Leve Item ID 12863)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(78, 8).Value = 12000
$ws.Cells.Item(78, 10).Value = 12000
$ws.Cells.Item(78, 12).Value = 108000
$ws.Cells.Item(78, 14).Value = -117984

# CUL row 80 (Leve Item ID 12890)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 8214.286
$ws.Cells.Item(80, 9).Value = 3500
$ws.Cells.Item(80, 10).Value = 9000
$ws.Cells.Item(80, 11).Value = 10500
$ws.Cells.Item(80, 12).Value = 27000
$ws.Cells.Item(80, 13).Value = -9564
$ws.Cells.Item(80, 14).Value = -28872

# CUL row 83 (Leve Item ID 12890)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(83, 8).Value = 8214.286
$ws.Cells.Item(83, 9).Value = 3500
$ws.Cells.Item(83, 10).Value = 9000
$ws.Cells.Item(83, 11).Value = 31500
$ws.Cells.Item(83, 12).Value = 81000
$ws.Cells.Item(83, 13).Value = -26820
$ws.Cells.Item(83, 14).Value = -90360

# CUL row 123 (Leve Item ID 36037)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(123, 8).Value = 4325
$ws.Cells.Item(123, 9).Value = 2933.3333
$ws.Cells.Item(123, 11).Value = 8799.999899999999
$ws.Cells.Item(123, 13).Value = -6349.999899999999

# LTW row 16 (Leve Item ID 5289)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1385.2354
$ws.Cells.Item(16, 9).Value = 1284.3125
$ws.Cells.Item(16, 10).Value = 3000
$ws.Cells.Item(16, 11).Value = 1284.3125
$ws.Cells.Item(16, 12).Value = 3000
$ws.Cells.Item(16, 13).Value = -1114.3125
$ws.Cells.Item(16, 14).Value = -3340

# LTW row 93 (Leve Item ID 19993)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 1534.6666
$ws.Cells.Item(93, 9).Value = 1377.7142
$ws.Cells.Item(93, 10).Value = 2633.3333
$ws.Cells.Item(93, 11).Value = 1377.7142
$ws.Cells.Item(93, 12).Value = 2633.3333
$ws.Cells.Item(93, 13).Value = -129.7141999999999
$ws.Cells.Item(93, 14).Value = -5129.3333

# WVR row 126 (Leve Item ID 36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1782.75
$ws.Cells.Item(126, 9).Value = 1186.1578
$ws.Cells.Item(126, 10).Value = 3042.2222
$ws.Cells.Item(126, 11).Value = 3558.4734
$ws.Cells.Item(126, 12).Value = 9126.6666
$ws.Cells.Item(126, 13).Value = -1088.4734
$ws.Cells.Item(126, 14).Value = -14066.6666
